# Auto-applied crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.261.44"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "'1.831.72"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'243.35"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'0.6167"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07361"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "'0.2907"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'23.29"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").Value = "'0.07650"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'1.838.70"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'0.6756"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'82.69"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'0.000008948"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'29.246.04"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "'2.087.01"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "'236.99"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "'12.53"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'7.398"
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'158.74"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'8.558"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "'0.1394"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'17.66"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "'1.496"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'0.05806"
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("D31").Value = "'1.235"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").Value = "'4.096"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "'4.106"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").Value = "'1.862"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "'0.7233"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("D37").Value = "'2.616"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "'2.860"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").Value = "'1.222.84"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'0.01764"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "'2.006.36"
$ws.Range("D45").Value = "'101.91"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'65.74"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'0.00000000120"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'0.5051"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "'9.218"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'0.4044"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'0.1166"
$ws.Range("E51").Value = "  +5.22%  "

# Rows 41 and 42 swapped coins (FraxShare <-> TrustWalletToken)
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9090"
$ws.Range("E41").Value = "  +1.48%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.208"
$ws.Range("E42").Value = "  -3.33%  "
